# ============================================================
# Update gh-pages output (commit 456a3b4)
# Workbook: 上海-漫展信息.xlsx
# ============================================================
$wb = $excel.ActiveWorkbook

# --- Sheet 1 (展览): refresh "想去人数" (F column) counts ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 274
$ws1.Range("F3").Value = 944
$ws1.Range("F4").Value = 567
$ws1.Range("F5").Value = 2327
$ws1.Range("F7").Value = 140
$ws1.Range("F8").Value = 1215
$ws1.Range("F10").Value = 3135
$ws1.Range("F13").Value = 1126
$ws1.Range("F14").Value = 650
$ws1.Range("F15").Value = 22
$ws1.Range("F17").Value = 262
$ws1.Range("F18").Value = 20
$ws1.Range("F20").Value = 1257
$ws1.Range("F21").Value = 1257
$ws1.Range("F22").Value = 201
$ws1.Range("F24").Value = 565
$ws1.Range("F25").Value = 235
$ws1.Range("F27").Value = 387
$ws1.Range("F28").Value = 669
$ws1.Range("F29").Value = 614
$ws1.Range("F30").Value = 265
$ws1.Range("F32").Value = 868
$ws1.Range("F33").Value = 142
$ws1.Range("F35").Value = 347
$ws1.Range("F36").Value = 1087
$ws1.Range("F37").Value = 5145
$ws1.Range("F38").Value = 589
$ws1.Range("F39").Value = 326
$ws1.Range("F40").Value = 208

# --- Sheet 2 (演出): remove two cancelled/merged listings (old rows 3 & 4), which
#     shifts every later listing up by two rows; the running index in column A is
#     then rewritten so it stays a plain 1..N sequence. ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A3:A4").EntireRow.Delete()
$lastRow = $ws2.Cells(($ws2.Rows.Count), 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 1
}

# --- Sheet 3 (本地生活): refresh "想去人数" (F column) counts ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 667
$ws3.Range("F5").Value = 469
$ws3.Range("F6").Value = 459

# --- Sheet 4 (全部类型): refresh "想去人数" (F column) counts ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 667
$ws4.Range("F3").Value = 274
$ws4.Range("F4").Value = 469
$ws4.Range("F6").Value = 944
$ws4.Range("F7").Value = 2327
$ws4.Range("F9").Value = 140
$ws4.Range("F10").Value = 1215
$ws4.Range("F13").Value = 3135
$ws4.Range("F17").Value = 650
$ws4.Range("F18").Value = 1781
$ws4.Range("F20").Value = 262
$ws4.Range("F22").Value = 20
$ws4.Range("F23").Value = 1257
$ws4.Range("F24").Value = 1257
$ws4.Range("F26").Value = 565
$ws4.Range("F27").Value = 312
$ws4.Range("F28").Value = 235
$ws4.Range("F29").Value = 387
$ws4.Range("F31").Value = 669
$ws4.Range("F32").Value = 614
$ws4.Range("F34").Value = 265
$ws4.Range("F35").Value = 868
$ws4.Range("F38").Value = 347
$ws4.Range("F39").Value = 1087
$ws4.Range("F40").Value = 68
$ws4.Range("F41").Value = 465
$ws4.Range("F42").Value = 326
$ws4.Range("F43").Value = 208
$ws4.Range("F47").Value = 31
$ws4.Range("F48").Value = 769
